$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Update "Story Points Completed" values for rows 5 and 6
$ws.Range("H5").Value = 1.5
$ws.Range("H6").Value = 1.5

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("H2:H6").Select()
